$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), styled like the existing headers (B1:H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$hdr = $ws.Range("I1:J1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# Data values for the new I and J columns, rows 2-30 (both columns hold the same value per row)
$values = @(
    7,
    8,
    8,
    7,
    8,
    9,
    10,
    10,
    8,
    9,
    9,
    9,
    10,
    9,
    9,
    9,
    9,
    9,
    8,
    9,
    9,
    9,
    9,
    9,
    7,
    6,
    8,
    5,
    6
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i]
    $ws.Cells.Item($row, 10).Value = $values[$i]
}
